# Add data for 2025-02-07
# Applies cell-level updates across Citywide Totals, By Neighborhood, and
# individual neighborhood sheets for the violent-crime-full-year workbook.

$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 577
$ws.Range("J3").Value = 8078
$ws.Range("L3").Value = 565
$ws.Range("H4").Value = 1748
$ws.Range("L4").Value = 154
$ws.Range("K6").Value = 9123
$ws.Range("L6").Value = 648
$ws.Range("H7").Value = 26060
$ws.Range("J7").Value = 29323
$ws.Range("K7").Value = 27527
$ws.Range("L7").Value = 1990

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 72
$ws.Range("L8").Value = 107
$ws.Range("K19").Value = 793
$ws.Range("L19").Value = 69
$ws.Range("L20").Value = 54
$ws.Range("L23").Value = 20
$ws.Range("L29").Value = 104
$ws.Range("L30").Value = 14
$ws.Range("L33").Value = 83
$ws.Range("L36").Value = 33
$ws.Range("K37").Value = 912
$ws.Range("L37").Value = 65
$ws.Range("L42").Value = 71
$ws.Range("L43").Value = 17
$ws.Range("L44").Value = 12
$ws.Range("L50").Value = 16
$ws.Range("L54").Value = 38
$ws.Range("L55").Value = 21
$ws.Range("L57").Value = 15
$ws.Range("L60").Value = 14
$ws.Range("J63").Value = 201
$ws.Range("K63").Value = 79
$ws.Range("L63").Value = 13
$ws.Range("L67").Value = 59
$ws.Range("L69").Value = 7
$ws.Range("L77").Value = 13
$ws.Range("L79").Value = 55
$ws.Range("L83").Value = 40
$ws.Range("H85").Value = 1231
$ws.Range("K85").Value = 1275
$ws.Range("L85").Value = 98
$ws.Range("L86").Value = 16
$ws.Range("L89").Value = 22
$ws.Range("L91").Value = 32
$ws.Range("L94").Value = 28
$ws.Range("K95").Value = 460
$ws.Range("L95").Value = 25
$ws.Range("L96").Value = 16
$ws.Range("L98").Value = 17
$ws.Range("L99").Value = 36
$ws.Range("H101").Value = 26060
$ws.Range("J101").Value = 29323
$ws.Range("K101").Value = 27527
$ws.Range("L101").Value = 1990

# --- West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 16

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 72

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 22

# --- South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 25
$ws.Range("H4").Value = 62
$ws.Range("K4").Value = 61
$ws.Range("H7").Value = 1231
$ws.Range("K7").Value = 1275
$ws.Range("L7").Value = 98

# --- Norwood Park ---
$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 7

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 29
$ws.Range("L4").Value = 8
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 107

# --- South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 40

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 83

# --- West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 9
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 460
$ws.Range("L7").Value = 25

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 20
$ws.Range("K4").Value = 42
$ws.Range("L4").Value = 4
$ws.Range("K7").Value = 912
$ws.Range("L7").Value = 65

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 36

# --- Fuller Park ---
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 14

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 59

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 11
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 38

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 104

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K4").Value = 33
$ws.Range("L6").Value = 24
$ws.Range("K7").Value = 793
$ws.Range("L7").Value = 69

# --- Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 12

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 71

# --- Lower West Side ---
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 21

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 20

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 12
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 32

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 55

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 54

# --- Grand Boulevard ---
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 33

# --- West Loop ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 28

# --- Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 17

# --- Lincoln Square ---
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 16

# --- Streeterville ---
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 16

# --- Mckinley Park ---
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 15

# --- Morgan Park ---
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 14

# --- Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 17

# --- Riverdale ---
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 13
